$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve existing formatting on the Price/Volume columns, force text
# entry so numeric-looking strings (e.g. "42.635.25", "0.636") are not
# auto-converted to numbers/dates by Excel, then restore original style.
$numRng = $ws.Range("D2:E51")
$origStyle = $numRng.Style
$numRng.NumberFormat = "@"

$ws.Range('D2').Value = '42.635.25'
$ws.Range('E2').Value = '  +1.06%  '

$ws.Range('D3').Value = '2.280.78'
$ws.Range('E3').Value = '  +2.57%  '

$ws.Range('E4').Value = '  -0.03%  '

$ws.Range('D5').Value = '251.49'
$ws.Range('E5').Value = '  -0.23%  '

$ws.Range('D6').Value = '0.636'
$ws.Range('E6').Value = '  +1.84%  '

$ws.Range('D7').Value = '73.74'
$ws.Range('E7').Value = '  +8.12%  '

$ws.Range('E8').Value = '  -0.20%  '

$ws.Range('D9').Value = '0.649'
$ws.Range('E9').Value = '  +3.30%  '

$ws.Range('D10').Value = '39.01'
$ws.Range('E10').Value = '  -0.14%  '

$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').Value = '0.0974'
$ws.Range('E11').Value = '  +3.44%  '

$ws.Range('B12').Value = 'OKB'
$ws.Range('C12').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D12').Value = '59.12'
$ws.Range('E12').Value = '  -1.24%  '

$ws.Range('D13').Value = '7.37'
$ws.Range('E13').Value = '  +3.78%  '

$ws.Range('E14').Value = '  +2.04%  '

$ws.Range('D15').Value = '2.624.33'
$ws.Range('E15').Value = '  +2.63%  '

$ws.Range('D16').Value = '15.04'
$ws.Range('E16').Value = '  +2.72%  '

$ws.Range('D17').Value = '0.878'
$ws.Range('E17').Value = '  -0.10%  '

$ws.Range('D18').Value = '2.296.99'
$ws.Range('E18').Value = '  +3.50%  '

$ws.Range('D19').Value = '42.580.32'
$ws.Range('E19').Value = '  +1.17%  '

$ws.Range('D20').Value = '0.0000100'
$ws.Range('E20').Value = '  +3.64%  '

$ws.Range('D21').Value = '6.30'
$ws.Range('E21').Value = '  +1.60%  '

$ws.Range('D22').Value = '72.06'
$ws.Range('E22').Value = '  -0.88%  '

$ws.Range('D23').Value = '235.53'
$ws.Range('E23').Value = '  +1.40%  '

$ws.Range('D24').Value = '2.21'
$ws.Range('E24').Value = '  +8.68%  '

$ws.Range('E25').Value = '  +0.30%  '

$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').Value = '11.52'
$ws.Range('E26').Value = '  +1.13%  '

$ws.Range('B27').Value = 'Dai'
$ws.Range('C27').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D27').Value = '0.998'
$ws.Range('E27').Value = '  -0.43%  '

$ws.Range('D28').Value = '2.44'
$ws.Range('E28').Value = '  +0.56%  '

$ws.Range('D29').Value = '3.66'

$ws.Range('D30').Value = '2.19'
$ws.Range('E30').Value = '  +6.30%  '

$ws.Range('D31').Value = '167.32'
$ws.Range('E31').Value = '  -0.08%  '

$ws.Range('D32').Value = '21.07'
$ws.Range('E32').Value = '  +2.89%  '

$ws.Range('D33').Value = '6.37'
$ws.Range('E33').Value = '  +7.27%  '

$ws.Range('D34').Value = '0.126'
$ws.Range('E34').Value = '  +4.08%  '

$ws.Range('D35').Value = '0.0805'
$ws.Range('E35').Value = '  +2.21%  '

$ws.Range('D36').Value = '31.54'
$ws.Range('E36').Value = '  +20.45%  '

$ws.Range('E37').Value = '  +3.02%  '

$ws.Range('D38').Value = '4.70'
$ws.Range('E38').Value = '  +13.65%  '

$ws.Range('D39').Value = '4.76'
$ws.Range('E39').Value = '  +2.77%  '

$ws.Range('D40').Value = '0.0307'
$ws.Range('E40').Value = '  -0.38%  '

$ws.Range('D41').Value = '13.66'
$ws.Range('E41').Value = '  +11.32%  '

$ws.Range('D42').Value = '2.35'
$ws.Range('E42').Value = '  +4.30%  '

$ws.Range('D43').Value = '5.93'
$ws.Range('E43').Value = '  +4.09%  '

$ws.Range('E44').Value = '  +7.16%  '

$ws.Range('D45').Value = '9.16'
$ws.Range('E45').Value = '  +6.45%  '

$ws.Range('D46').Value = '61.92'
$ws.Range('E46').Value = '  -0.30%  '

$ws.Range('E47').Value = '  -3.66%  '

$ws.Range('E48').Value = '  +3.64%  '

$ws.Range('E49').Value = '  -0.12%  '

$ws.Range('E50').Value = '  +1.23%  '

$ws.Range('D51').Value = '97.77'
$ws.Range('E51').Value = '  +4.39%  '

# Restore original cell formatting
$numRng.Style = $origStyle
